$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns remain stored as text (matches source data which is
# inline text, not numeric), so values like "1.00" or "261.20" keep their exact
# textual representation instead of being coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '37.667.53'
$ws.Range("E2").Value = '  +2.83%  '
$ws.Range("D3").Value = '2.057.50'
$ws.Range("E3").Value = '  +4.30%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '261.20'
$ws.Range("E5").Value = '  +6.91%  '
$ws.Range("E6").Value = '  -0.28%  '
$ws.Range("D7").Value = '58.79'
$ws.Range("E7").Value = '  -2.45%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("E9").Value = '  +3.41%  '
$ws.Range("D10").Value = '57.44'
$ws.Range("E10").Value = '  +0.85%  '
$ws.Range("D11").Value = '0.0813'
$ws.Range("E11").Value = '  +3.12%  '
$ws.Range("E12").Value = '  -0.04%  '
$ws.Range("D13").Value = '15.16'
$ws.Range("E13").Value = '  +6.56%  '
$ws.Range("D14").Value = '2.350.17'
$ws.Range("E14").Value = '  +3.84%  '
$ws.Range("E15").Value = '  -0.12%  '
$ws.Range("D16").Value = '21.74'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("E17").Value = '  +0.33%  '
$ws.Range("D18").Value = '2.047.74'
$ws.Range("E18").Value = '  +3.80%  '
$ws.Range("D19").Value = '37.638.34'
$ws.Range("E19").Value = '  +3.02%  '
$ws.Range("D20").Value = '70.65'
$ws.Range("E20").Value = '  +1.29%  '
$ws.Range("D21").Value = '0.0₃0864'
$ws.Range("E21").Value = '  +1.12%  '
$ws.Range("D22").Value = '5.33'
$ws.Range("E22").Value = '  +4.58%  '
$ws.Range("D23").Value = '230.30'
$ws.Range("E23").Value = '  +0.34%  '
$ws.Range("D24").Value = '2.68'
$ws.Range("E24").Value = '  +9.72%  '
$ws.Range("D25").Value = '1.01'
$ws.Range("E25").Value = '  +0.54%  '
$ws.Range("D26").Value = '2.36'
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("D27").Value = '9.32'
$ws.Range("E27").Value = '  +1.35%  '
$ws.Range("E28").Value = '  -2.98%  '
$ws.Range("D29").Value = '164.49'
$ws.Range("E29").Value = '  +1.81%  '
$ws.Range("D30").Value = '20.09'
$ws.Range("E30").Value = '  +3.72%  '
$ws.Range("E31").Value = '  +1.62%  '
$ws.Range("D32").Value = '0.121'
$ws.Range("E32").Value = '  +0.77%  '
$ws.Range("E33").Value = '  +1.36%  '
$ws.Range("D34").Value = '0.0675'
$ws.Range("E34").Value = '  +9.65%  '
$ws.Range("D35").Value = '4.57'
$ws.Range("E35").Value = '  +1.47%  '
$ws.Range("D36").Value = '2.53'
$ws.Range("E36").Value = '  +11.14%  '
$ws.Range("D37").Value = '3.55'
$ws.Range("E37").Value = '  +7.58%  '
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.05%  '
$ws.Range("E39").Value = '  +2.27%  '
$ws.Range("E40").Value = '  +0.93%  '
$ws.Range("E41").Value = '  +4.01%  '
$ws.Range("D42").Value = '0.0982'
$ws.Range("E42").Value = '  +0.78%  '
$ws.Range("D43").Value = '0.0219'
$ws.Range("E43").Value = '  +4.68%  '
$ws.Range("D44").Value = '1.20'
$ws.Range("E44").Value = '  +3.26%  '
$ws.Range("D45").Value = '16.78'
$ws.Range("E45").Value = '  +5.97%  '
$ws.Range("D46").Value = '1.410.45'
$ws.Range("E46").Value = '  +3.25%  '
$ws.Range("D47").Value = '92.71'
$ws.Range("E47").Value = '  +4.15%  '
$ws.Range("E48").Value = '  +4.31%  '
$ws.Range("D49").Value = '7.54'
$ws.Range("E49").Value = '  +4.94%  '
$ws.Range("E50").Value = '  +11.28%  '
$ws.Range("E51").Value = '  +2.90%  '
